$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column G: "supplier_name" header + "Global Suppliers" for every data row ---
$ws.Range("G1").Value = "supplier_name"
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "Global Suppliers"
}

# --- Styles for the new column (use a temporary named style so the engine does not
#     fork an extra "unlinked from theme" clone of the default font) ---
$headerStyle = $wb.Styles.Add("TempSupplierHeaderStyle")
$headerStyle.Font.Name = "Segoe UI"
$headerStyle.Font.Size = 11
$ws.Range("G1").Style = "TempSupplierHeaderStyle"
$headerStyle.Delete()

$dataStyle = $wb.Styles.Add("TempSupplierDataStyle")
$dataStyle.Font.Name = "Segoe UI"
$dataStyle.Font.Size = 8
$ws.Range("G2:G51").Style = "TempSupplierDataStyle"
$dataStyle.Delete()

# --- Column width for G ---
$ws.Columns.Item(7).ColumnWidth = 13.166666666666666

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 16.5
$ws.Range("A2:A51").EntireRow.RowHeight = 14.65

# --- Selection / scroll position ---
$ws.Range("G1").Select() | Out-Null
